$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.431.89"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "2.954.94"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.00"
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.88"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "2.951.34"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("E10").Value = "  -4.30%  "
$ws.Range("E11").Value = "  -3.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000244"
$ws.Range("E13").Value = "  -2.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.60"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").Value = "65.371.12"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").Value = "3.441.53"
$ws.Range("E17").Value = "  -1.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.04"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.02"
$ws.Range("E19").Value = "  +15.30%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "2.949.82"
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.22"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.698"
$ws.Range("E22").Value = "  +2.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.31"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.45"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.22"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.01"
$ws.Range("E28").Value = "  -5.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.48"
$ws.Range("E29").Value = "  +6.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.96"
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.59"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.113"
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.19"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.976"
$ws.Range("E36").Value = "  -1.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.74"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "44.77"
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("E40").Value = "  -7.35%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.121"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.301"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.84"
$ws.Range("E43").Value = "  -4.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.56"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "383.38"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0351"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").Value = "2.683.65"
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.46"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.68"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.18"
$ws.Range("E51").Value = "  +1.79%  "